# Portfolio Optimisation / annual_returns.xlsx - "monte carlo and portfolio
# optimisation last trys"
#
# Adds a new "Sheet3" (Date / IIH / IPB) built from Sheet2's Date, IIH and
# IPB columns (i.e. Sheet2 minus its TLE column), makes it the active sheet
# with the selection on G11, and leaves Sheet2 with a "select all" selection
# instead of its previous single-cell selection.

$wb = $excel.ActiveWorkbook
$sheet2 = $wb.Worksheets.Item("Sheet2")

# --- Create the new sheet right after Sheet2 ---------------------------
$sheet3 = $wb.Worksheets.Add($null, $sheet2)
$sheet3.Name = "Sheet3"

# --- Header row ----------------------------------------------------------
$sheet3.Range("A1").Value = "Date"
$sheet3.Range("A1").NumberFormat = "[$-409]mmm\-yy;@"
$sheet3.Range("A1").Font.Bold = $true
$sheet3.Range("A1").HorizontalAlignment = -4108

$sheet3.Range("B1").Value = "IIH"
$sheet3.Range("B1").Font.Bold = $true
$sheet3.Range("B1").HorizontalAlignment = -4108

$sheet3.Range("C1").Value = "IPB"
$sheet3.Range("C1").Font.Bold = $true
$sheet3.Range("C1").HorizontalAlignment = -4108

# --- Data rows (Date, IIH, IPB) taken from Sheet2's A/C/D columns -------
$data = @(
    @(43831.0, 4.68, 0.621),
    @(43862.0, 1.24, 0.103),
    @(43891.0, 3.82, -0.079),
    @(43922.0, 13.69, 1.292),
    @(43952.0, 8.091, -1.072),
    @(43983.0, 10.228, 0.361),
    @(44013.0, -3.285, 3.317),
    @(44044.0, -2.709, 3.462),
    @(44075.0, 4.759, 11.193),
    @(44105.0, 1.451, 8.383),
    @(44136.0, 13.269, 1.886),
    @(44166.0, 20.662, 7.019),
    @(44197.0, 4.538, 4.852),
    @(44228.0, 1.58, 3.859),
    @(44256.0, -2.448, 1.826),
    @(44287.0, 2.126, 1.511),
    @(44317.0, 4.98, 6.334),
    @(44348.0, -5.443, -1.453),
    @(44378.0, 3.923, 0.074),
    @(44409.0, 5.634, 1.019),
    @(44440.0, -3.541, 4.29),
    @(44470.0, 7.761, 2.4),
    @(44501.0, 16.68, 19.061),
    @(44531.0, 3.878, -8.491),
    @(44562.0, 9.408, 7.084),
    @(44593.0, -4.383, -7.301),
    @(44621.0, 14.97, 14.288),
    @(44652.0, 12.061, 12.873),
    @(44682.0, 3.202, 4.353),
    @(44713.0, -3.802, -0.763),
    @(44743.0, 10.0, 8.416),
    @(44774.0, 24.41, 16.191),
    @(44805.0, 0.745, 9.314),
    @(44835.0, 24.412, 18.811),
    @(44866.0, 26.701, 23.09),
    @(44896.0, 12.461, 16.107),
    @(44927.0, -9.627, -9.475),
    @(44958.0, 1.336, 3.02),
    @(44986.0, -1.491, 4.112),
    @(45017.0, 2.829, 4.871),
    @(45047.0, 4.831, 7.768),
    @(45078.0, 23.011, 23.794),
    @(45108.0, 28.764, 23.854),
    @(45139.0, 22.725, 4.607),
    @(45170.0, 8.559, 14.397),
    @(45200.0, -5.175, -1.941),
    @(45231.0, 6.239, 11.045),
    @(45261.0, -1.316, -0.378),
    @(45292.0, 14.401, 4.916),
    @(45323.0, 6.436, 10.709),
    @(45352.0, 4.266, 3.79),
    @(45383.0, 16.942, 7.355),
    @(45413.0, 8.78, 4.127),
    @(45444.0, 7.055, 5.255),
    @(45474.0, 0.127, 7.559),
)

$r = 2
foreach ($row in $data) {
    $sheet3.Cells.Item($r, 1).Value = $row[0]
    $sheet3.Cells.Item($r, 1).NumberFormat = "[$-409]mmm\-yy;@"
    $sheet3.Cells.Item($r, 1).Font.Bold = $true

    $sheet3.Cells.Item($r, 2).Value = $row[1]
    $sheet3.Cells.Item($r, 3).Value = $row[2]
    $r++
}

# --- Selections: Sheet2 becomes "select all", Sheet3 is the new active
#     sheet with G11 selected -------------------------------------------
$sheet2.Activate()
$sheet2.Cells.Select()

$sheet3.Activate()
$sheet3.Range("G11").Select()
